# Add "Wins" / "Losses" / "Ties" columns (AD, AE, AF) to the sheet.
# Column header row (row 1) gets the same bold/bordered/centered style as
# the other header cells (e.g. style index 1 in the original workbook:
# bold font, thin box border, horizontal=center, vertical=top).
# Data rows 2-46 get the season record values: Wins=81, Losses=81, Ties=0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# --- Data rows: every team played an 81-81-0 season record ---
for ($row = 2; $row -le 46; $row++) {
    $ws.Cells.Item($row, 30).Value = 81
    $ws.Cells.Item($row, 31).Value = 81
    $ws.Cells.Item($row, 32).Value = 0
}

Write-Output "Added Wins/Losses/Ties columns (AD:AF) for rows 1-46"
